$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"   = 8.807400000000003
    "B7"   = 4.9678
    "B16"  = 6.779200000000002
    "B28"  = 5.961599999999997
    "B29"  = 5.114600000000001
    "B32"  = 7.569199999999992
    "B40"  = 8.875799999999996
    "B52"  = 5.030400000000003
    "B57"  = 5.182499999999997
    "B66"  = 5.385899999999998
    "B100" = 5.316799999999998
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
